$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185, shifting the existing rows (185-269) down to (186-270)
$ws.Rows.Item(185).Insert()

# Populate the newly inserted row 185 with this week's new record
$ws.Range("A185").Value = 5
$ws.Range("B185").Value = "Macroferia Regional de Talca"
$ws.Range("C185").Value = "Maule"
$ws.Range("D185").Value = 45016
$ws.Range("E185").Value = 7
$ws.Range("F185").Value = 100112017
$ws.Range("G185").Value = "Apio"
$ws.Range("H185").Value = "Americana (o)"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 500
$ws.Range("K185").Value = 8000
$ws.Range("L185").Value = 8000
$ws.Range("M185").Value = 8000
$ws.Range("N185").Value = "`$/docena de matas"
$ws.Range("O185").Value = "Provincia del Elquí"
$ws.Range("P185").Value = 1333
$ws.Range("Q185").Value = 6
$ws.Range("R185").Value = "Hortaliza"
